# "auto AGV push away"
# Renumber Tank1/Tank2 end-times and bump AGV1's start index, then
# adjust the corresponding numeric rows so the job data stays consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings) to reflect new timing.
$ws.Range("C1").Value = "Tank1|3"
$ws.Range("D1").Value = "Tank2|4"
$ws.Range("H1").Value = "AGV1|2"

# Job 1 row: AGV1 assignment moves from slot 3 to slot 1.
$ws.Range("B2").Value = 1

# Job 1 row: Tank1/Tank2 duration shrinks from 20 to 10.
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 10

# Job 2 row: AGV1 assignment moves from slot 3 to slot 1.
$ws.Range("B6").Value = 1

# Job 2 row: Tank1/Tank2 duration shrinks from 40 to 12.
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 12
